{"js": "// Remove the literal \" {location}\" placeholder text (a space run followed by\n// the \"{location}\" run) that used to sit between \"{company}\" and the\n// tab character in the work-history paragraph. After the edit the line\n// reads \"{company}\\t{startDate} - {endDate}\".\nconst body = context.document.body;\nconst results = body.search(\" {location}\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the literal \" {location}\" placeholder text (a space followed by the\n# \"{location}\" token) that used to sit between \"{company}\" and the tab\n# character in the work-history paragraph. After the edit the line reads\n# \"{company}`t{startDate} - {endDate}\".\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \" {location}\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Wrap = 1\n\n$found = $find.Execute()\nif ($found) {\n    $rng.Delete()\n}\n"}
